# Updated cryptos list — refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active sheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "26.668.98";  E = "  -0.45%  " },
    @{ Row = 3;  D = "1.597.10";   E = "  -0.70%  " },
    @{ Row = 4;  D = $null;        E = "  +0.02%  " },
    @{ Row = 5;  D = "211.27";     E = "  +0.26%  " },
    @{ Row = 6;  D = "0.511";      E = "  +0.44%  " },
    @{ Row = 7;  D = $null;        E = "  -0.03%  " },
    @{ Row = 8;  D = "0.0618";     E = "  -0.45%  " },
    @{ Row = 9;  D = $null;        E = "  -1.25%  " },
    @{ Row = 10; D = "19.71";      E = "  +0.00%  " },
    @{ Row = 11; D = $null;        E = "  +0.11%  " },
    @{ Row = 12; D = "1.820.82";   E = $null },
    @{ Row = 13; D = "1.617.67";   E = "  +0.45%  " },
    @{ Row = 14; D = $null;        E = "  -0.78%  " },
    @{ Row = 15; D = $null;        E = "  -1.80%  " },
    @{ Row = 16; D = "65.10";      E = "  +2.10%  " },
    @{ Row = 17; D = "26.650.04";  E = "  -0.51%  " },
    @{ Row = 18; D = "0.0₃0727";   E = "  -0.23%  " },
    @{ Row = 19; D = "209.64";     E = "  -0.23%  " },
    @{ Row = 20; D = $null;        E = "  +0.03%  " },
    @{ Row = 21; D = "6.77";       E = "  +0.39%  " },
    @{ Row = 22; D = $null;        E = "  -0.49%  " },
    @{ Row = 23; D = $null;        E = "  -1.26%  " },
    @{ Row = 24; D = "8.92";       E = "  +0.55%  " },
    @{ Row = 25; D = "146.48";     E = "  -0.05%  " },
    @{ Row = 26; D = $null;        E = "  +0.17%  " },
    @{ Row = 27; D = "7.17";       E = "  -4.06%  " },
    @{ Row = 28; D = $null;        E = "  +2.17%  " },
    @{ Row = 29; D = "15.30";      E = "  -0.27%  " },
    @{ Row = 30; D = $null;        E = "  +0.16%  " },
    @{ Row = 31; D = $null;        E = "  +0.04%  " },
    @{ Row = 32; D = $null;        E = "  -1.13%  " },
    @{ Row = 33; D = "0.666";      E = "  -2.20%  " },
    @{ Row = 34; D = $null;        E = "  -1.34%  " },
    @{ Row = 35; D = "1.297.96";   E = "  -1.61%  " },
    @{ Row = 36; D = "2.44";       E = "  +0.08%  " },
    @{ Row = 37; D = $null;        E = "  -3.09%  " },
    @{ Row = 38; D = $null;        E = "  -1.40%  " },
    @{ Row = 39; D = "0.843";      E = "  +2.48%  " },
    @{ Row = 40; D = $null;        E = "  +0.04%  " },
    @{ Row = 41; D = $null;        E = "  +2.03%  " },
    @{ Row = 42; D = $null;        E = "  +0.30%  " },
    @{ Row = 43; D = "0.789";      E = "  -0.25%  " },
    @{ Row = 44; D = "63.76";      E = "  +0.89%  " },
    @{ Row = 45; D = "1.733.58";   E = "  -0.84%  " },
    @{ Row = 46; D = "90.22";      E = "  +1.28%  " },
    @{ Row = 47; D = "0.880";      E = "  +7.82%  " },
    @{ Row = 48; D = $null;        E = "  +0.96%  " },
    @{ Row = 49; D = $null;        E = "  +2.96%  " },
    @{ Row = 50; D = $null;        E = "  -1.40%  " },
    @{ Row = 51; D = "7.48";       E = "  -0.16%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
